$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.039913953218831
$ws.Range("D2").Value = 1.043047995616468
$ws.Range("E2").Value = 1.047803087651592
$ws.Range("F2").Value = 1.057294117141863
$ws.Range("I2").Value = 1.040920471246013
$ws.Range("J2").Value = 1.045003393771235
$ws.Range("K2").Value = 1.045822961556264
$ws.Range("L2").Value = 1.050564705429315
$ws.Range("M2").Value = 1.060029481313317
$ws.Range("N2").Value = 1.005712725503983
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.04077876476748
$ws.Range("D3").Value = 1.043702153493407
$ws.Range("E3").Value = 1.048583400728539
$ws.Range("F3").Value = 1.058181755642709
$ws.Range("I3").Value = 1.041126990835344
$ws.Range("J3").Value = 1.045514030462325
$ws.Range("K3").Value = 1.046288485527669
$ws.Range("L3").Value = 1.051157018058541
$ws.Range("M3").Value = 1.060730740038175
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.041338989996178
$ws.Range("D4").Value = 1.044125928433514
$ws.Range("E4").Value = 1.049089277680736
$ws.Range("F4").Value = 1.058757269054845
$ws.Range("I4").Value = 1.041259771337901
$ws.Range("J4").Value = 1.045844419810791
$ws.Range("K4").Value = 1.046589515701228
$ws.Range("L4").Value = 1.051540585455487
$ws.Range("M4").Value = 1.061185013231135
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.041574659011278
$ws.Range("D5").Value = 1.044304199282138
$ws.Range("E5").Value = 1.049302176952785
$ws.Range("F5").Value = 1.058999488836887
$ws.Range("I5").Value = 1.04131538772072
$ws.Range("J5").Value = 1.045983307862026
$ws.Range("K5").Value = 1.04671602087977
$ws.Range("L5").Value = 1.05170190798808
$ws.Range("M5").Value = 1.061376110671841
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.041614237646861
$ws.Range("D6").Value = 1.044334138487063
$ws.Range("E6").Value = 1.049337937031653
$ws.Range("F6").Value = 1.059040174614979
$ws.Range("I6").Value = 1.041324713941632
$ws.Range("J6").Value = 1.046006627293747
$ws.Range("K6").Value = 1.046737258810854
$ws.Range("L6").Value = 1.05172899886258
$ws.Range("M6").Value = 1.061408203822587
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.041342138428311
$ws.Range("D7").Value = 1.044128310044182
$ws.Range("E7").Value = 1.049092121555183
$ws.Range("F7").Value = 1.058760504533315
$ws.Range("I7").Value = 1.041260515291297
$ws.Range("J7").Value = 1.045846275671253
$ws.Range("K7").Value = 1.046591206259398
$ws.Range("L7").Value = 1.051542740778608
$ws.Range("M7").Value = 1.061187566211036
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.040206088067412
$ws.Range("D8").Value = 1.043268968639539
$ws.Range("E8").Value = 1.048066597821806
$ws.Range("F8").Value = 1.057593859107094
$ws.Range("I8").Value = 1.040990441395057
$ws.Range("J8").Value = 1.045175970627664
$ws.Range("K8").Value = 1.045980327190048
$ws.Range("L8").Value = 1.05076481699089
$ws.Range("M8").Value = 1.060266368050175
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.03820914816888
$ws.Range("D9").Value = 1.041758535082429
$ws.Range("E9").Value = 1.046266940983305
$ws.Range("F9").Value = 1.055546984556841
$ws.Range("I9").Value = 1.04050804611637
$ws.Range("J9").Value = 1.043994657162587
$ws.Range("K9").Value = 1.044902441503179
$ws.Range("L9").Value = 1.049396386447454
$ws.Range("M9").Value = 1.058647088860734
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.036881257933301
$ws.Range("D10").Value = 1.04075425943529
$ws.Range("E10").Value = 1.045072281360775
$ws.Range("F10").Value = 1.054188497365829
$ws.Range("I10").Value = 1.040182126833608
$ws.Range("J10").Value = 1.043207087088128
$ws.Range("K10").Value = 1.044182957158901
$ws.Range("L10").Value = 1.048485775847303
$ws.Range("M10").Value = 1.057570344927497
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.036307093053398
$ws.Range("D11").Value = 1.040320054014509
$ws.Range("E11").Value = 1.044556215049024
$ws.Range("F11").Value = 1.053601726352139
$ws.Range("I11").Value = 1.04003998368607
$ws.Range("J11").Value = 1.042866069399387
$ws.Range("K11").Value = 1.043871215947937
$ws.Range("L11").Value = 1.048091886438726
$ws.Range("M11").Value = 1.057104780858035
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.036093947425671
$ws.Range("D12").Value = 1.040158870554561
$ws.Range("E12").Value = 1.044364711502949
$ws.Range("F12").Value = 1.053383995123005
$ws.Range("I12").Value = 1.039987033144472
$ws.Range("J12").Value = 1.042739402285815
$ws.Range("K12").Value = 1.043755392612211
$ws.Range("L12").Value = 1.047945641521638
$ws.Range("M12").Value = 1.056931952348123
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.036139662215581
$ws.Range("D13").Value = 1.040193440411531
$ws.Range("E13").Value = 1.044405781197239
$ws.Range("F13").Value = 1.053430689150589
$ws.Range("I13").Value = 1.039998398093858
$ws.Range("J13").Value = 1.042766572699494
$ws.Range("K13").Value = 1.04378023839288
$ws.Range("L13").Value = 1.047977008672756
$ws.Range("M13").Value = 1.056969019976736
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.036289471815148
$ws.Range("D14").Value = 1.040306728494247
$ws.Range("E14").Value = 1.044540381492097
$ws.Range("F14").Value = 1.053583724088398
$ws.Range("I14").Value = 1.040035609879606
$ws.Range("J14").Value = 1.0428555990075
$ws.Range("K14").Value = 1.043861642534598
$ws.Range("L14").Value = 1.048079796491972
$ws.Range("M14").Value = 1.057090492694999
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.036381791033139
$ws.Range("D15").Value = 1.040376542292607
$ws.Range("E15").Value = 1.044623337937372
$ws.Range("F15").Value = 1.053678043391178
$ws.Range("I15").Value = 1.040058517135257
$ws.Range("J15").Value = 1.042910451361392
$ws.Range("K15").Value = 1.043911794537805
$ws.Range("L15").Value = 1.048143135869338
$ws.Range("M15").Value = 1.057165349703814
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.036919380779358
$ws.Range("D16").Value = 1.040783090115158
$ws.Range("E16").Value = 1.0451065570325
$ws.Range("F16").Value = 1.05422747040653
$ws.Range("I16").Value = 1.040191539017816
$ws.Range("J16").Value = 1.043229719527221
$ws.Range("K16").Value = 1.044203642293809
$ws.Range("L16").Value = 1.048511925768732
$ws.Range("M16").Value = 1.057601257199028
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.037256817331171
$ws.Range("D17").Value = 1.041038282751759
$ws.Range("E17").Value = 1.045409997955483
$ws.Range("F17").Value = 1.054572504435851
$ws.Range("I17").Value = 1.040274708081472
$ws.Range("J17").Value = 1.043429990317349
$ws.Range("K17").Value = 1.044386658078625
$ws.Range("L17").Value = 1.048743368951554
$ws.Range("M17").Value = 1.057874871793119
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.037453717412751
$ws.Range("D18").Value = 1.041187195146208
$ws.Range("E18").Value = 1.045587108390976
$ws.Range("F18").Value = 1.054773897971482
$ws.Range("I18").Value = 1.040323120925809
$ws.Range("J18").Value = 1.043546805314227
$ws.Range("K18").Value = 1.044493388734164
$ws.Range("L18").Value = 1.048878405405032
$ws.Range("M18").Value = 1.058034531400699
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.037520868614138
$ws.Range("D19").Value = 1.041237981031883
$ws.Range("E19").Value = 1.045647518514209
$ws.Range("F19").Value = 1.054842591825985
$ws.Range("I19").Value = 1.040339611766574
$ws.Range("J19").Value = 1.043586636245305
$ws.Range("K19").Value = 1.04452977780884
$ws.Range("L19").Value = 1.048924456013589
$ws.Range("M19").Value = 1.058088982151531
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.037220605395125
$ws.Range("D20").Value = 1.04101089648395
$ws.Range("E20").Value = 1.045377429348448
$ws.Range("F20").Value = 1.05453547094965
$ws.Range("I20").Value = 1.040265794986932
$ws.Range("J20").Value = 1.043408503103173
$ws.Range("K20").Value = 1.044367024210032
$ws.Range("L20").Value = 1.048718533186343
$ws.Range("M20").Value = 1.05784550881811
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.036245353165363
$ws.Range("D21").Value = 1.040273365209623
$ws.Range("E21").Value = 1.044500739898518
$ws.Range("F21").Value = 1.053538652979852
$ws.Range("I21").Value = 1.040024656133433
$ws.Range("J21").Value = 1.042829382912474
$ws.Range("K21").Value = 1.043837671828035
$ws.Range("L21").Value = 1.04804952628671
$ws.Range("M21").Value = 1.057054719158973
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.035632896038937
$ws.Range("D22").Value = 1.039810228095669
$ws.Range("E22").Value = 1.04395061013062
$ws.Range("F22").Value = 1.052913197289821
$ws.Range("I22").Value = 1.039872162049297
$ws.Range("J22").Value = 1.04246527947795
$ws.Range("K22").Value = 1.04350468113699
$ws.Range("L22").Value = 1.047629261140087
$ws.Range("M22").Value = 1.05655811337287
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.035957502139895
$ws.Range("D23").Value = 1.040055690472105
$ws.Range("E23").Value = 1.044242141390608
$ws.Range("F23").Value = 1.053244640916547
$ws.Range("I23").Value = 1.039953085281196
$ws.Range("J23").Value = 1.042658296007872
$ws.Range("K23").Value = 1.043681221108257
$ws.Range("L23").Value = 1.047852016481387
$ws.Range("M23").Value = 1.056821316563448
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.037236967766223
$ws.Range("D24").Value = 1.041023270965578
$ws.Range("E24").Value = 1.045392145334213
$ws.Range("F24").Value = 1.054552204352468
$ws.Range("I24").Value = 1.040269822734394
$ws.Range("J24").Value = 1.043418212249163
$ws.Range("K24").Value = 1.044375895970075
$ws.Range("L24").Value = 1.048729755277534
$ws.Range("M24").Value = 1.057858776482439
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.03872481187598
$ws.Range("D25").Value = 1.042148553380382
$ws.Range("E25").Value = 1.046731302121333
$ws.Range("F25").Value = 1.05607508427776
$ws.Range("I25").Value = 1.040633521760836
$ws.Range("J25").Value = 1.044300065275101
$ws.Range("K25").Value = 1.045181263722695
$ws.Range("L25").Value = 1.04974986911252
$ws.Range("M25").Value = 1.059065228961354
